# Update the GitHub-link shape on slide 9 ("Content Placeholder 2"):
#  - nudge the shape's horizontal position by 1 EMU (834504 -> 834503 EMU)
#  - widen the shape (8915400 -> 10832777 EMU)
#  - bump the run's font size (20pt -> 32pt)
#  - fix/extend the URL text (lowercase "marcos" -> "Marcos", add ".git")
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)

# Position / size are exposed by the COM object model in points
# (1 pt = 12700 EMU). Use high-precision point values so the
# underlying EMU values come out to exactly 834503 / 10832777.
$shp.Left = 65.70889763779527
$shp.Width = 852.9746

$tr = $shp.TextFrame.TextRange
$tr.Text = "https://github.com/Surajmohapatra49/Marcos-project.git"
$tr.Font.Size = 32
